$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (number formats / styles) from the last existing data row (187)
# down onto the two new rows (188, 189) so they reuse the workbook's existing
# cell styles (date style on column A, number style on columns C:Z) instead of
# Excel fabricating brand-new style/numFmt entries.
$ws.Range("A187:Z187").Copy()
$ws.Range("A188:Z189").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 188: 2025-12-03, 四方坪站充电量(kw)
$ws.Range("A188").Value = 45994
$ws.Range("B188").Value = "四方坪站充电量(kw)"
$ws.Range("C188").Value = 481.87
$ws.Range("D188").Value = 915.02100000000019
$ws.Range("E188").Value = 503.20000000000005
$ws.Range("F188").Value = 327.57
$ws.Range("G188").Value = 417.61600000000004
$ws.Range("H188").Value = 530.22299999999996
$ws.Range("I188").Value = 315.45199999999994
$ws.Range("J188").Value = 83.567999999999998
$ws.Range("K188").Value = 90.38
$ws.Range("L188").Value = 183.59
$ws.Range("M188").Value = 261.18299999999999
$ws.Range("N188").Value = 167.66199999999998
$ws.Range("O188").Value = 1072.4230000000002
$ws.Range("P188").Value = 1237.7700000000004
$ws.Range("Q188").Value = 449.3830000000001
$ws.Range("R188").Value = 378.47800000000001
$ws.Range("S188").Value = 419.28099999999984
$ws.Range("T188").Value = 289.17899999999997
$ws.Range("U188").Value = 154.28
$ws.Range("V188").Value = 86.08
$ws.Range("W188").Value = 7.84
$ws.Range("X188").Value = 139.35500000000002
$ws.Range("Y188").Value = 88.09
$ws.Range("Z188").Value = 53.72

# Row 189: 2025-12-03, 高岭站充电量(kw)
$ws.Range("A189").Value = 45994
$ws.Range("B189").Value = "高岭站充电量(kw)"
$ws.Range("C189").Value = 372.02700000000004
$ws.Range("D189").Value = 437.34699999999998
$ws.Range("E189").Value = 105.28800000000001
$ws.Range("F189").Value = 0
$ws.Range("G189").Value = 97.366
$ws.Range("H189").Value = 67.149000000000001
$ws.Range("I189").Value = 133.25800000000001
$ws.Range("J189").Value = 130.70599999999999
$ws.Range("K189").Value = 325.29100000000005
$ws.Range("L189").Value = 283.97199999999998
$ws.Range("M189").Value = 158.89600000000002
$ws.Range("N189").Value = 285.10100000000006
$ws.Range("O189").Value = 588.39800000000002
$ws.Range("P189").Value = 458.75500000000005
$ws.Range("Q189").Value = 395.76900000000001
$ws.Range("R189").Value = 201.839
$ws.Range("S189").Value = 87.503999999999991
$ws.Range("T189").Value = 99.073000000000008
$ws.Range("U189").Value = 0
$ws.Range("V189").Value = 0
$ws.Range("W189").Value = 28.992000000000001
$ws.Range("X189").Value = 53.403999999999996
$ws.Range("Y189").Value = 0
$ws.Range("Z189").Value = 42.36

# Update the active selection to match the post-edit state (E193)
$ws.Range("E193").Select()
